$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.524.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.958.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("E6").Value = '  +1.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.59'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.76%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +4.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0787'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.103'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.21'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.838'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.246.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.36'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.959.78'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.494.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0848'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '229.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.53%  '
$ws.Range("E24").Value = '  +5.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.50%  '
$ws.Range("E26").Value = '  +8.09%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.120'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.72%  '
$ws.Range("E31").Value = '  +8.70%  '
$ws.Range("E32").Value = '  +4.39%  '
$ws.Range("E33").Value = '  -1.56%  '
$ws.Range("E34").Value = '  +6.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +19.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.64%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("E38").Value = '  -1.24%  '
$ws.Range("E39").Value = '  -9.37%  '
$ws.Range("E40").Value = '  +1.54%  '
$ws.Range("E41").Value = '  +1.13%  '
$ws.Range("E42").Value = '  +2.44%  '
$ws.Range("E43").Value = '  +1.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.370.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.75'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.18%  '
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("E48").Value = '  +1.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.84'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.137.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.95'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.66%  '
